# Insert a new daily-stats row right before the 2026/12/29 block
# (row 685) for date 2026/01/24 (土), time 3, ranking 17.
# Every existing row from 685 downward shifts down by one (685->686 ... 726->727),
# and the sheet's used dimension grows from D726 to D727.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 685..726 down to 686..727, creating a blank row 685.
$ws.Rows.Item(685).Insert()

# The date column stores plain text like "2026/12/29"; format the cell as
# Text first so Excel doesn't auto-convert the literal into a date serial.
$ws.Range("A685").NumberFormat = "@"
$ws.Range("A685").Value = "2026/01/24"
$ws.Range("B685").Value = "土"
$ws.Range("C685").Value = 3
$ws.Range("D685").Value = 17

# Re-align the new cell's style with its neighbors so it doesn't keep a
# stray "Text" number format / quote-prefix style that the surrounding
# (unstyled) data cells don't have.
$ws.Range("A685").Style = $ws.Range("A686").Style
$ws.Range("B685").Style = $ws.Range("B686").Style
